$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2").Value = 0.09702050429457308
$ws.Range("J2").Value = 0.09731161031231482
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3387653333333334
$ws.Range("N2").Value = 1.016296
$ws.Range("O2").Value = 0.1020856984836854
$ws.Range("P2").Value = 0.1038222745268702
$ws.Range("Q2").Value = 0.074212982808
$ws.Range("R2").Value = 0.667916845272
$ws.Range("S2").Value = 0.009904405948150893
$ws.Range("T2").Value = 0.01010311272049697
$ws.Range("I3").Value = 0.09702050429457308
$ws.Range("J3").Value = 0.09731161031231482
$ws.Range("M3").Value = 1.622048333333334
$ws.Range("N3").Value = 4.866145
$ws.Range("O3").Value = 0.4887983532827969
$ws.Range("P3").Value = 0.4971132840014691
$ws.Range("Q3").Value = 0.355340506335
$ws.Range("R3").Value = 3.198064557015
$ws.Range("S3").Value = 0.04742346273385385
$ws.Range("T3").Value = 0.04837489417382605
$ws.Range("I4").Value = 0.09702050429457308
$ws.Range("J4").Value = 0.09731161031231482
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6269413333333334
$ws.Range("N4").Value = 1.880824
$ws.Range("O4").Value = 0.1889264857530475
$ws.Range("P4").Value = 0.1921403072182968
$ws.Range("Q4").Value = 0.137343410952
$ws.Range("R4").Value = 1.236090698568
$ws.Range("S4").Value = 0.01832974292236214
$ws.Range("T4").Value = 0.01869748270131535
$ws.Range("I5").Value = 0.09702050429457308
$ws.Range("J5").Value = 0.09731161031231482
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.166517
$ws.Range("N5").Value = 0.333034
$ws.Range("O5").Value = 0.05017929103649284
$ws.Range("P5").Value = 0.03402192606758435
$ws.Range("Q5").Value = 0.036478712673
$ws.Range("R5").Value = 0.218872276038
$ws.Range("S5").Value = 0.004868420121504686
$ws.Range("T5").Value = 0.003310728411563154
$ws.Range("I6").Value = 0.09702050429457308
$ws.Range("J6").Value = 0.09731161031231482
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5641686666666667
$ws.Range("N6").Value = 1.692506
$ws.Range("O6").Value = 0.1700101714439774
$ws.Range("P6").Value = 0.1729022081857796
$ws.Range("Q6").Value = 0.123591865638
$ws.Range("R6").Value = 1.112326790742
$ws.Range("S6").Value = 0.01649447256870151
$ws.Range("T6").Value = 0.01682539230511331
$ws.Range("G7").Value = 1.165950333333333
$ws.Range("H7").Value = 3.497851
$ws.Range("I7").Value = 0.516371961904357
$ws.Range("J7").Value = 0.5179213146581529
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3387653333333334
$ws.Range("N7").Value = 1.016296
$ws.Range("O7").Value = 0.1020856984836854
$ws.Range("P7").Value = 0.1038222745268702
$ws.Range("Q7").Value = 0.3949835533217778
$ws.Range("R7").Value = 3.554851979896
$ws.Range("S7").Value = 0.05271419240839727
$ws.Range("T7").Value = 0.0537717689137563
$ws.Range("G8").Value = 1.165950333333333
$ws.Range("H8").Value = 3.497851
$ws.Range("I8").Value = 0.516371961904357
$ws.Range("J8").Value = 0.5179213146581529
$ws.Range("M8").Value = 1.622048333333334
$ws.Range("N8").Value = 4.866145
$ws.Range("O8").Value = 0.4887983532827969
$ws.Range("P8").Value = 0.4971132840014691
$ws.Range("Q8").Value = 1.891227794932778
$ws.Range("R8").Value = 17.021050154395
$ws.Range("S8").Value = 0.2524017646602568
$ws.Range("T8").Value = 0.2574655655840726
$ws.Range("G9").Value = 1.165950333333333
$ws.Range("H9").Value = 3.497851
$ws.Range("I9").Value = 0.516371961904357
$ws.Range("J9").Value = 0.5179213146581529
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6269413333333334
$ws.Range("N9").Value = 1.880824
$ws.Range("O9").Value = 0.1889264857530475
$ws.Range("P9").Value = 0.1921403072182968
$ws.Range("Q9").Value = 0.7309824565804444
$ws.Range("R9").Value = 6.578842109223999
$ws.Range("S9").Value = 0.09755634010399666
$ws.Range("T9").Value = 0.09951356051332169
$ws.Range("G10").Value = 1.165950333333333
$ws.Range("H10").Value = 3.497851
$ws.Range("I10").Value = 0.516371961904357
$ws.Range("J10").Value = 0.5179213146581529
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.166517
$ws.Range("N10").Value = 0.333034
$ws.Range("O10").Value = 0.05017929103649284
$ws.Range("P10").Value = 0.03402192606758435
$ws.Range("Q10").Value = 0.1941505516556666
$ws.Range("R10").Value = 1.164903309934
$ws.Range("S10").Value = 0.02591117895948352
$ws.Range("T10").Value = 0.01762068067612577
$ws.Range("G11").Value = 1.165950333333333
$ws.Range("H11").Value = 3.497851
$ws.Range("I11").Value = 0.516371961904357
$ws.Range("J11").Value = 0.5179213146581529
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.5641686666666667
$ws.Range("N11").Value = 1.692506
$ws.Range("O11").Value = 0.1700101714439774
$ws.Range("P11").Value = 0.1729022081857796
$ws.Range("Q11").Value = 0.6577926449562221
$ws.Range("R11").Value = 5.920133804606
$ws.Range("S11").Value = 0.08778848577222267
$ws.Range("T11").Value = 0.08954973897087662
$ws.Range("G12").Value = 0.020264
$ws.Range("H12").Value = 0.040528
$ws.Range("I12").Value = 0.008974448685232641
$ws.Range("J12").Value = 0.006000917432007717
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.3387653333333334
$ws.Range("N12").Value = 1.016296
$ws.Range("O12").Value = 0.1020856984836854
$ws.Range("P12").Value = 0.1038222745268702
$ws.Range("Q12").Value = 0.006864740714666667
$ws.Range("R12").Value = 0.041188444288
$ws.Range("S12").Value = 0.0009161628625379664
$ws.Range("T12").Value = 0.0006230288970389864
$ws.Range("G13").Value = 0.020264
$ws.Range("H13").Value = 0.040528
$ws.Range("I13").Value = 0.008974448685232641
$ws.Range("J13").Value = 0.006000917432007717
$ws.Range("M13").Value = 1.622048333333334
$ws.Range("N13").Value = 4.866145
$ws.Range("O13").Value = 0.4887983532827969
$ws.Range("P13").Value = 0.4971132840014691
$ws.Range("Q13").Value = 0.03286918742666668
$ws.Range("R13").Value = 0.19721512456
$ws.Range("S13").Value = 0.004386695738962676
$ws.Range("T13").Value = 0.002983135771647018
$ws.Range("G14").Value = 0.020264
$ws.Range("H14").Value = 0.040528
$ws.Range("I14").Value = 0.008974448685232641
$ws.Range("J14").Value = 0.006000917432007717
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.6269413333333334
$ws.Range("N14").Value = 1.880824
$ws.Range("O14").Value = 0.1889264857530475
$ws.Range("P14").Value = 0.1921403072182968
$ws.Range("Q14").Value = 0.01270433917866667
$ws.Range("R14").Value = 0.076226035072
$ws.Range("S14").Value = 0.00169551105167206
$ws.Range("T14").Value = 0.001153018118977596
$ws.Range("G15").Value = 0.020264
$ws.Range("H15").Value = 0.040528
$ws.Range("I15").Value = 0.008974448685232641
$ws.Range("J15").Value = 0.006000917432007717
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.5
$ws.Range("M15").Value = 0.166517
$ws.Range("N15").Value = 0.333034
$ws.Range("O15").Value = 0.05017929103649284
$ws.Range("P15").Value = 0.03402192606758435
$ws.Range("Q15").Value = 0.003374300488
$ws.Range("R15").Value = 0.013497201952
$ws.Range("S15").Value = 0.0004503314724683592
$ws.Range("T15").Value = 0.0002041627692094447
$ws.Range("G16").Value = 0.020264
$ws.Range("H16").Value = 0.040528
$ws.Range("I16").Value = 0.008974448685232641
$ws.Range("J16").Value = 0.006000917432007717
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.5641686666666667
$ws.Range("N16").Value = 1.692506
$ws.Range("O16").Value = 0.1700101714439774
$ws.Range("P16").Value = 0.1729022081857796
$ws.Range("Q16").Value = 0.01143231386133333
$ws.Range("R16").Value = 0.06859388316800001
$ws.Range("S16").Value = 0.001525747559591579
$ws.Range("T16").Value = 0.001037571875134672
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8526826666666666
$ws.Range("H17").Value = 2.558048
$ws.Range("I17").Value = 0.3776330851158373
$ws.Range("J17").Value = 0.3787661575975245
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.3387653333333334
$ws.Range("N17").Value = 1.016296
$ws.Range("O17").Value = 0.1020856984836854
$ws.Range("P17").Value = 0.1038222745268702
$ws.Range("Q17").Value = 0.2888593278008889
$ws.Range("R17").Value = 2.599733950208
$ws.Range("S17").Value = 0.03855093726459928
$ws.Range("T17").Value = 0.03932436399557799
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.8526826666666666
$ws.Range("H18").Value = 2.558048
$ws.Range("I18").Value = 0.3776330851158373
$ws.Range("J18").Value = 0.3787661575975245
$ws.Range("M18").Value = 1.622048333333334
$ws.Range("N18").Value = 4.866145
$ws.Range("O18").Value = 0.4887983532827969
$ws.Range("P18").Value = 0.4971132840014691
$ws.Range("Q18").Value = 1.383092498328889
$ws.Range("R18").Value = 12.44783248496
$ws.Range("S18").Value = 0.1845864301497235
$ws.Range("T18").Value = 0.1882896884719234
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.8526826666666666
$ws.Range("H19").Value = 2.558048
$ws.Range("I19").Value = 0.3776330851158373
$ws.Range("J19").Value = 0.3787661575975245
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.6269413333333334
$ws.Range("N19").Value = 1.880824
$ws.Range("O19").Value = 0.1889264857530475
$ws.Range("P19").Value = 0.1921403072182968
$ws.Range("Q19").Value = 0.5345820079502221
$ws.Range("R19").Value = 4.811238071552
$ws.Range("S19").Value = 0.07134489167501661
$ws.Range("T19").Value = 0.07277624588468219
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.8526826666666666
$ws.Range("H20").Value = 2.558048
$ws.Range("I20").Value = 0.3776330851158373
$ws.Range("J20").Value = 0.3787661575975245
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.5
$ws.Range("M20").Value = 0.166517
$ws.Range("N20").Value = 0.333034
$ws.Range("O20").Value = 0.05017929103649284
$ws.Range("P20").Value = 0.03402192606758435
$ws.Range("Q20").Value = 0.1419861596053333
$ws.Range("R20").Value = 0.8519169576319999
$ws.Range("S20").Value = 0.01894936048303627
$ws.Range("T20").Value = 0.01288635421068598
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.8526826666666666
$ws.Range("H21").Value = 2.558048
$ws.Range("I21").Value = 0.3776330851158373
$ws.Range("J21").Value = 0.3787661575975245
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.5641686666666667
$ws.Range("N21").Value = 1.692506
$ws.Range("O21").Value = 0.1700101714439774
$ws.Range("P21").Value = 0.1729022081857796
$ws.Range("Q21").Value = 0.481056843143111
$ws.Range("R21").Value = 4.329511588288
$ws.Range("T21").Value = 0.06548950503465499
